# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" right before the "总计" (Total) sheet,
#    containing the per-fund holding detail for the new quarter.
# 2. Insert a new summary row at the top of the "总计" sheet's data
#    (row 2) for "2022-Q1", shifting the existing quarters down and
#    renumbering the index column.

$wb = $excel.ActiveWorkbook
$total = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------
# Step 1: new "2022-Q1" sheet (fund-level detail), placed right before 总计
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Add($total)
$q1.Name = "2022-Q1"

$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $cell = $q1.Cells.Item(1, $c + 2)
    $cell.Value = $headers[$c]
    $cell.Font.Bold = $true
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

# code, name, scale, stockPosition, positionPct, marketValue, rank
$rows = @(
    @("001113","南方大数据100指数A","20.79","94.23","1.53","0.3181",10),
    @("001907","国投瑞银境煊灵活配置混合A","2.61","90.44","4.15","0.1083",10),
    @("001908","国投瑞银境煊灵活配置混合C","1.75","90.44","4.15","0.0726",10),
    @("015309","国投瑞银境煊灵活配置混合E","0.33","90.44","4.15","0.0137",10),
    @("004344","南方大数据100指数C","0.17","94.23","1.53","0.0026",10),
    @("004794","富荣福鑫灵活配置混合A","0.06","89.60","2.92","0.0018",8),
    @("004795","富荣福鑫灵活配置混合C","0.06","89.60","2.92","0.0018",8)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $row = $rows[$i]

    $idxCell = $q1.Cells.Item($r, 1)
    $idxCell.Value = $i
    $idxCell.Font.Bold = $true
    $idxCell.Borders.LineStyle = 1
    $idxCell.HorizontalAlignment = -4108
    $idxCell.VerticalAlignment = -4160

    $codeCell = $q1.Cells.Item($r, 2)
    $codeCell.NumberFormat = "@"
    $codeCell.Value = $row[0]

    $q1.Cells.Item($r, 3).Value = $row[1]

    for ($c = 4; $c -le 7; $c++) {
        $textCell = $q1.Cells.Item($r, $c)
        $textCell.NumberFormat = "@"
        $textCell.Value = $row[$c - 2]
    }

    $q1.Cells.Item($r, 8).Value = $row[6]
}

# ---------------------------------------------------------------------
# Step 2: insert the "2022-Q1" summary row into "总计"
# ---------------------------------------------------------------------
# Re-fetch the handle: the $total reference captured before the sheet
# insertion above can end up pointing at the newly added sheet instead.
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()
$total.Range("A2:D2").ClearFormats()

$total.Cells.Item(2, 2).Value = "2022-Q1"
$total.Cells.Item(2, 3).Value = 7
$total.Cells.Item(2, 4).Value = 0.52

for ($r = 2; $r -le 7; $r++) {
    $total.Cells.Item($r, 1).Value = $r - 2
}

$a2 = $total.Cells.Item(2, 1)
$a2.Font.Bold = $true
$a2.Borders.LineStyle = 1
$a2.HorizontalAlignment = -4108
$a2.VerticalAlignment = -4160

Write-Host "done"
